$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values.
# The NumberFormat/Style dance keeps the underlying cell style unchanged (s="0")
# while still storing the value as a string, matching the source data shape.
$textForceCells = @(
    'D11', 'D14', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D31', 'D32', 'D33', 'D36', 'D42', 'D45', 'D46', 'D47', 'D48', 'D49', 'D5', 'D6'
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '70.933.80'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '3.848.43'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '707.06'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').Value = '172.48'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '3.845.84'
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').Value = '7.31'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '36.69'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '4.495.37'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '3.863.77'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').Value = '70.981.87'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').Value = '17.35'
$ws.Range('E20').Value = '  -2.89%  '
$ws.Range('D21').Value = '10.67'
$ws.Range('E21').Value = '  -3.70%  '
$ws.Range('D22').Value = '493.13'
$ws.Range('E22').Value = '  +2.30%  '
$ws.Range('D23').Value = '0.716'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = '85.29'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('D26').Value = '10.64'
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').Value = '  +2.01%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '7.49'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('D33').Value = '29.46'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('D35').Value = '3.803.49'
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('D36').Value = '9.14'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +7.15%  '
$ws.Range('E40').Value = '  +6.82%  '
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('D42').Value = '3.32'
$ws.Range('E42').Value = '  -3.31%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = '0.000315'
$ws.Range('E45').Value = '  -3.60%  '
$ws.Range('D46').Value = '162.86'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '48.72'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '416.44'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '1.39'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('E51').Value = '  +0.63%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
